{"js": "// Update the date in the title paragraph.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items[0].insertText(\"2023-12-14 Thursday\", \"Replace\");\nawait context.sync();\n\n// Update the division-problem answers in the table, cell by cell\n// (position-based, so identical/overlapping old & new values never collide).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// row index -> new text for the 5 cells in that row (0-based row/col).\nconst rowUpdates = {\n  0: [\"71\u00f78=8, 7\", \"64\u00f77=9, 1\", \"83\u00f75=16, 3\", \"88\u00f73=29, 1\", \"52\u00f75=10, 2\"],\n  4: [\"86\u00f76=14, 2\", \"35\u00f75=7, 0\", \"66\u00f73=22, 0\", \"27\u00f72=13, 1\", \"12\u00f72=6, 0\"],\n  8: [\"77\u00f76=12, 5\", \"88\u00f75=17, 3\", \"48\u00f75=9, 3\", \"50\u00f72=25, 0\", \"55\u00f73=18, 1\"],\n  12: [\"75\u00f73=25, 0\", \"15\u00f74=3, 3\", \"99\u00f74=24, 3\", \"76\u00f79=8, 4\", \"35\u00f73=11, 2\"],\n  16: [\"26\u00f79=2, 8\", \"46\u00f76=7, 4\", \"60\u00f79=6, 6\", \"20\u00f76=3, 2\", \"87\u00f79=9, 6\"],\n};\n\nfor (const rowIndex of Object.keys(rowUpdates)) {\n  const r = Number(rowIndex);\n  const values = rowUpdates[rowIndex];\n  for (let c = 0; c < values.length; c++) {\n    table.getCell(r, c).value = values[c];\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date in the title paragraph.\n$d = $word.ActiveDocument\n$d.Paragraphs(1).Range.Text = \"2023-12-14 Thursday\"\n\n# Update the division-problem answers in the table, cell by cell\n# (position-based via Table.Cell(row, col), 1-based, so identical/\n# overlapping old & new values never collide).\n$t = $d.Tables(1)\n\n$rowUpdates = @{\n    1  = @(\"71\u00f78=8, 7\", \"64\u00f77=9, 1\", \"83\u00f75=16, 3\", \"88\u00f73=29, 1\", \"52\u00f75=10, 2\")\n    5  = @(\"86\u00f76=14, 2\", \"35\u00f75=7, 0\", \"66\u00f73=22, 0\", \"27\u00f72=13, 1\", \"12\u00f72=6, 0\")\n    9  = @(\"77\u00f76=12, 5\", \"88\u00f75=17, 3\", \"48\u00f75=9, 3\", \"50\u00f72=25, 0\", \"55\u00f73=18, 1\")\n    13 = @(\"75\u00f73=25, 0\", \"15\u00f74=3, 3\", \"99\u00f74=24, 3\", \"76\u00f79=8, 4\", \"35\u00f73=11, 2\")\n    17 = @(\"26\u00f79=2, 8\", \"46\u00f76=7, 4\", \"60\u00f79=6, 6\", \"20\u00f76=3, 2\", \"87\u00f79=9, 6\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $values = $rowUpdates[$rowIndex]\n    for ($c = 0; $c -lt $values.Count; $c++) {\n        $t.Cell($rowIndex, $c + 1).Range.Text = $values[$c]\n    }\n}\n"}
